$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy style from H1 header cell)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
    @(7, 7),
    @(9, 9),
    @(6, 8),
    @(6, 6),
    @(5, 5),
    @(6, 7),
    @(6, 8),
    @(9, 9),
    @(7, 8),
    @(5, 6),
    @(5, 6),
    @(6, 7),
    @(6, 7),
    @(9, 9),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
